# critical-path-method.xlsx — "Modification chemin critique et Planification sprints"
#
# Adds a new "Complétion" summary task (id 400) with its single sub-task
# "GP33 - Intégration" (id 410, predecessors 100/200/300, O/M/P = 4/6/8h),
# rewires the Finish milestone to depend on the new task 400 instead of
# directly on 100/200/300, swaps the Tests group's predecessor order,
# hard-codes the Finish Date cell, and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CPM")

# --- "Tests" group (row 18): predecessors 220,210 -> 210,220 ------------
$ws.Range("C18").Value = 210
$ws.Range("D18").Value = 220

# --- New row 26: id 400 "Complétion" (group header, predecessor = 410) --
$ws.Range("A26").Value = 400
$ws.Range("B26").Value = "Complétion"
$ws.Range("C26").Value = 410

# --- New row 27: id 410 "GP33 - Intégration" -----------------------------
#     predecessors 100, 200, 300 ; Optimistic/Most likely/Pessimistic = 4/6/8
$ws.Range("A27").Value = 410
$ws.Range("B27").Value = "GP33 - Intégration"
$ws.Range("C27").Value = 100
$ws.Range("D27").Value = 200
$ws.Range("E27").Value = 300
$ws.Range("I27").Value = 4
$ws.Range("J27").Value = 6
$ws.Range("K27").Value = 8

# --- Finish (row 33): predecessor 100,200,300 -> just 400 ---------------
$ws.Range("C33").Value = 400
$ws.Range("D33").Value = ""
$ws.Range("E33").Value = ""

# --- Finish Date (B6) is hard-coded instead of computed by WORKDAY() ----
$ws.Range("B6").Value = 44980

# --- Update the active selection / scroll position -----------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D19").Select()
